# ng-content - directive to display data in other components
# Add two new vocabulary rows (11 and 12) below the existing table in
# columns C/D, mirroring the formatting already used for row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary pairs
$ws.Range("C11").Value = "efficient"
$ws.Range("D11").Value = "wydajny"
$ws.Range("C12").Value = "get dumped"
$ws.Range("D12").Value = "zostać porzuconym"

# Copy the existing D10 formatting (wrap-text style) down onto D11:D12
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D11:D12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the active selection to the last edited cell, like the author did
$ws.Range("D12").Select() | Out-Null
